# Insert a new data row at row 441 (pushing existing rows 441..520 down to 442..521)
# and populate it with the new "Poroto granado" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(441).Insert()

$ws.Cells.Item(441, 1).Value = 6
$ws.Cells.Item(441, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(441, 3).Value = "Metropolitana"
$ws.Cells.Item(441, 4).Value = 44694
$ws.Cells.Item(441, 5).Value = 13
$ws.Cells.Item(441, 6).Value = 100112030
$ws.Cells.Item(441, 7).Value = "Poroto granado"
$ws.Cells.Item(441, 8).Value = "Sin especificar"
$ws.Cells.Item(441, 9).Value = "Primera"
$ws.Cells.Item(441, 10).Value = 290
$ws.Cells.Item(441, 11).Value = 20000
$ws.Cells.Item(441, 12).Value = 22000
$ws.Cells.Item(441, 13).Value = 20828
$ws.Cells.Item(441, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(441, 15).Value = "Región Metropolitana"
$ws.Cells.Item(441, 16).Value = 833
$ws.Cells.Item(441, 17).Value = 25
$ws.Cells.Item(441, 18).Value = "Hortaliza"
